$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53..236 down to 54..237.
$ws.Rows(53).Insert()

# Populate the newly inserted row 53 with the new record's data.
$ws.Range("A53").Value = 3
$ws.Range("B53").Value = "Femacal de La Calera"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 44560
$ws.Range("E53").Value = 5
$ws.Range("F53").Value = 100112039
$ws.Range("G53").Value = "Ciboulette"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 180
$ws.Range("K53").Value = 1500
$ws.Range("L53").Value = 1500
$ws.Range("M53").Value = 1500
$ws.Range("N53").Value = "$/docena de atados"
$ws.Range("O53").Value = "Provincia de Quillota"
$ws.Range("P53").Value = 500
$ws.Range("Q53").Value = 3
$ws.Range("R53").Value = "Hortaliza"
